# Updated cryptos list on Wed May  3 07:56:33 UTC 2023 with GitHub Actions
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row.
# NumberFormat is forced to Text ("@") before writing so that numeric-looking
# strings (e.g. "7.000", "0.9811") are preserved verbatim instead of being
# auto-coerced into numbers/percentages by Excel's input parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = "28.730.81";   E = "  +2.36%  " },
    @{ Row = 3;  D = "1.871.08";    E = "  +1.87%  " },
    @{ Row = 4;  D = $null;         E = "  +0.19%  " },
    @{ Row = 5;  D = "324.73";      E = "  +0.07%  " },
    @{ Row = 6;  D = "1.003";       E = "  +0.19%  " },
    @{ Row = 7;  D = "0.4603";      E = "  -0.76%  " },
    @{ Row = 8;  D = "0.3862";      E = "  -0.28%  " },
    @{ Row = 9;  D = "0.07859";     E = "  -0.07%  " },
    @{ Row = 10; D = "0.9811";      E = "  +1.93%  " },
    @{ Row = 11; D = "21.82";       E = "  -0.49%  " },
    @{ Row = 12; D = "1.845.48";    E = "  +1.12%  " },
    @{ Row = 13; D = "7.000";       E = "  +1.01%  " },
    @{ Row = 14; D = "5.705";       E = "  +0.09%  " },
    @{ Row = 15; D = "0.06964";     E = "  +1.58%  " },
    @{ Row = 16; D = "88.42";       E = "  +0.93%  " },
    @{ Row = 17; D = "1.004";       E = "  +0.28%  " },
    @{ Row = 18; D = "0.00001002";  E = "  +0.70%  " },
    @{ Row = 19; D = "16.77";       E = "  +0.45%  " },
    @{ Row = 20; D = "1.003";       E = "  +0.17%  " },
    @{ Row = 21; D = "28.740.66";   E = "  +2.35%  " },
    @{ Row = 22; D = "5.272";       E = "  -1.19%  " },
    @{ Row = 23; D = "11.08";       E = "  +0.50%  " },
    @{ Row = 24; D = "2.099";       E = "  +0.06%  " },
    @{ Row = 25; D = "2.104.49";    E = "  +2.83%  " },
    @{ Row = 26; D = "152.45";      E = "  -1.30%  " },
    @{ Row = 27; D = "19.29";       E = "  +0.61%  " },
    @{ Row = 28; D = "5.891";       E = "  +3.64%  " },
    @{ Row = 29; D = "1.988";       E = "  +1.13%  " },
    @{ Row = 30; D = "119.15";      E = "  +0.72%  " },
    @{ Row = 31; D = "0.09325";     E = "  +0.95%  " },
    @{ Row = 32; D = "0.9163";      E = "  -2.47%  " },
    @{ Row = 33; D = "5.297";       E = "  +0.23%  " },
    @{ Row = 34; D = "1.333";       E = "  +0.72%  " },
    @{ Row = 35; D = "3.324";       E = "  +0.49%  " },
    @{ Row = 36; D = "0.05783";     E = "  -1.36%  " },
    @{ Row = 37; D = "1.146";       E = "  +0.62%  " },
    @{ Row = 38; D = "0.02079";     E = $null },
    @{ Row = 39; D = "7.653";       E = "  -1.84%  " },
    @{ Row = 40; D = "0.5626";      E = "  +0.49%  " },
    @{ Row = 41; D = "0.1779";      E = "  +0.82%  " },
    @{ Row = 42; D = "9.773";       E = "  -1.45%  " },
    @{ Row = 43; D = "0.07222";     E = "  -0.61%  " },
    @{ Row = 44; D = "11.72";       E = "  +0.83%  " },
    @{ Row = 45; D = "0.5291";      E = "  +0.35%  " },
    @{ Row = 46; D = "2.127";       E = "  +0.10%  " },
    @{ Row = 47; D = "1.123";       E = "  +0.23%  " },
    @{ Row = 48; D = $null;         E = "  +0.15%  " },
    @{ Row = 49; D = "112.92";      E = "  +0.25%  " },
    @{ Row = 50; D = $null;         E = "  +3.65%  " },
    @{ Row = 51; D = "1.003";       E = "  +0.22%  " }
)

foreach ($r in $rows) {
    if ($null -ne $r.D) {
        $cell = $ws.Range("D" + $r.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $r.D
    }
    if ($null -ne $r.E) {
        $cell = $ws.Range("E" + $r.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $r.E
    }
}
